$p = $ppt.ActivePresentation

# --- Slide 1: title shape, fix "Estimación SpaO2 y " -> "Estimación SpO2 y " ---
# Keep the whole phrase as a single run (same run boundaries as before), just
# correct the misspelling inside it.
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(2)
$titleRange = $titleShape.TextFrame.TextRange
$found1 = $titleRange.Find("Estimación SpaO2 y ")
if ($found1 -ne $null) {
    $found1.Text = "Estimación SpO2 y "
}

# --- Slide 7: two result tables, fix "Spao2 ..." header cells -> "SpO2 ..." ---
$s7 = $p.Slides.Item(7)
for ($i = 1; $i -le $s7.Shapes.Count; $i++) {
    $shp = $s7.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.Cell(1,2).Shape.TextFrame.TextRange.Text = "SpO2 Medido [%]"
        $tbl.Cell(1,3).Shape.TextFrame.TextRange.Text = "SpO2 Real [%]"
        $tbl.Cell(1,4).Shape.TextFrame.TextRange.Text = "SpO2 Error [%]"
    }
}

# --- Slide 8: conclusions text, split "Los programas de spaO2 y BR..." into
#     three runs so "SpO2" can carry its own (English) run ---
$s8 = $p.Slides.Item(8)
$bodyShape = $s8.Shapes.Item(3)
$bodyRange = $bodyShape.TextFrame.TextRange
$found2 = $bodyRange.Find("Los programas de spaO2 y BR son distintos, lo que disminuye los tiempos de ejecución (a 5[s] y 20[s] respectivamente).")
if ($found2 -ne $null) {
    # Normalise the whole sentence first (keeps it as a single run).
    $found2.Text = "Los programas de SpO2 y BR son distintos, lo que disminuye los tiempos de ejecución (a 5[s] y 20[s] respectivamente)."

    # Re-find "SpO2" inside the (now corrected) sentence and give it its own run.
    $spo2Range = $bodyRange.Find("SpO2")
    if ($spo2Range -ne $null) {
        $spo2Range.Text = "SpO2"
        $spo2Range.LanguageID = 1033
    }
}
